$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-16 04:19:03"
$ws.Range("E3").Value = "2026-02-16 04:19:05"
$ws.Range("E4").Value = "2026-02-16 04:19:08"
$ws.Range("H4").NumberFormat = "@"
$ws.Range("H4").Value = "68%"
$ws.Range("J4").Value = "1014.6 hPa"
$ws.Range("E5").Value = "2026-02-16 04:19:10"
$ws.Range("I5").Value = "1.7 mm"
$ws.Range("N5").Value = "-1.2 °C 3:58 TU"
$ws.Range("E6").Value = "2026-02-16 04:19:13"
$ws.Range("J6").Value = "1014.7 hPa"
$ws.Range("E7").Value = "2026-02-16 04:19:15"
$ws.Range("J7").Value = "1015.0 hPa"
$ws.Range("E8").Value = "2026-02-16 04:19:18"
$ws.Range("J8").Value = "1014.8 hPa"
$ws.Range("M8").Value = "9.6 °C 3:59 TU"
$ws.Range("O8").Value = "9.3 °C"
$ws.Range("E9").Value = "2026-02-16 04:19:21"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "95%"
$ws.Range("N9").Value = "3.8 °C 3:53 TU"
$ws.Range("O9").Value = "5.3 °C"
$ws.Range("E10").Value = "2026-02-16 04:19:23"
$ws.Range("E11").Value = "2026-02-16 04:19:26"
$ws.Range("N11").Value = "0.6 °C 3:58 TU"
$ws.Range("O11").Value = "0.9 °C"
$ws.Range("E12").Value = "2026-02-16 04:19:29"
$ws.Range("O12").Value = "5.6 °C"
$ws.Range("E13").Value = "2026-02-16 04:19:32"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "87%"
$ws.Range("O13").Value = "1.3 °C"
$ws.Range("E14").Value = "2026-02-16 04:19:34"
$ws.Range("E15").Value = "2026-02-16 04:19:37"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "89%"
$ws.Range("N15").Value = "4.0 °C 3:59 TU"
$ws.Range("O15").Value = "5.6 °C"
$ws.Range("E16").Value = "2026-02-16 04:19:40"
$ws.Range("M16").Value = "0.2 °C 3:58 TU"
$ws.Range("O16").Value = "-0.8 °C"
$ws.Range("E17").Value = "2026-02-16 04:19:43"
$ws.Range("E18").Value = "2026-02-16 04:19:45"
$ws.Range("J18").Value = "1015.1 hPa"
$ws.Range("O18").Value = "4.4 °C"
$ws.Range("E19").Value = "2026-02-16 04:19:48"
$ws.Range("H19").NumberFormat = "@"
$ws.Range("H19").Value = "95%"
$ws.Range("N19").Value = "2.9 °C 3:59 TU"
$ws.Range("E20").Value = "2026-02-16 04:19:51"
$ws.Range("H20").NumberFormat = "@"
$ws.Range("H20").Value = "90%"
$ws.Range("N20").Value = "-1.7 °C 3:49 TU"
$ws.Range("O20").Value = "-1.1 °C"
$ws.Range("E21").Value = "2026-02-16 04:19:53"
$ws.Range("O21").Value = "4.9 °C"
$ws.Range("E22").Value = "2026-02-16 04:19:56"
$ws.Range("N22").Value = "-6.5 °C 3:38 TU"
$ws.Range("E23").Value = "2026-02-16 04:19:59"
$ws.Range("O23").Value = "-1.0 °C"
$ws.Range("E24").Value = "2026-02-16 04:20:02"
$ws.Range("J24").Value = "1018.2 hPa"
$ws.Range("E25").Value = "2026-02-16 04:20:04"
$ws.Range("H25").NumberFormat = "@"
$ws.Range("H25").Value = "71%"
$ws.Range("O25").Value = "0.7 °C"
$ws.Range("E26").Value = "2026-02-16 04:20:07"
$ws.Range("E27").Value = "2026-02-16 04:20:09"
$ws.Range("N27").Value = "0.4 °C 3:57 TU"
$ws.Range("E28").Value = "2026-02-16 04:20:12"
$ws.Range("J28").Value = "1016.0 hPa"
$ws.Range("N28").Value = "2.6 °C 3:59 TU"
$ws.Range("E29").Value = "2026-02-16 04:20:15"
$ws.Range("O29").Value = "5.0 °C"
$ws.Range("E30").Value = "2026-02-16 04:20:17"
$ws.Range("J30").Value = "1014.7 hPa"
$ws.Range("E31").Value = "2026-02-16 04:20:20"
$ws.Range("J31").Value = "1013.2 hPa"
$ws.Range("O31").Value = "14.1 °C"
$ws.Range("E32").Value = "2026-02-16 04:20:23"
$ws.Range("E33").Value = "2026-02-16 04:20:26"
$ws.Range("H33").NumberFormat = "@"
$ws.Range("H33").Value = "71%"
$ws.Range("N33").Value = "2.8 °C 3:47 TU"
$ws.Range("O33").Value = "4.9 °C"
$ws.Range("E34").Value = "2026-02-16 04:20:29"
$ws.Range("H34").NumberFormat = "@"
$ws.Range("H34").Value = "64%"
$ws.Range("E35").Value = "2026-02-16 04:20:32"
$ws.Range("H35").NumberFormat = "@"
$ws.Range("H35").Value = "77%"
$ws.Range("J35").Value = "1019.2 hPa"
$ws.Range("N35").Value = "6.6 °C 3:41 TU"
$ws.Range("O35").Value = "6.8 °C"
$ws.Range("E36").Value = "2026-02-16 04:20:34"
$ws.Range("H36").NumberFormat = "@"
$ws.Range("H36").Value = "90%"
$ws.Range("J36").Value = "1014.6 hPa"
$ws.Range("O36").Value = "7.0 °C"
$ws.Range("E37").Value = "2026-02-16 04:20:37"
$ws.Range("N37").Value = "0.7 °C 3:55 TU"
$ws.Range("O37").Value = "1.9 °C"
$ws.Range("E38").Value = "2026-02-16 04:20:40"
$ws.Range("E39").Value = "2026-02-16 04:20:42"
$ws.Range("H39").NumberFormat = "@"
$ws.Range("H39").Value = "73%"
$ws.Range("M39").Value = "0.5 °C 3:32 TU"
$ws.Range("O39").Value = "-0.1 °C"
$ws.Range("E40").Value = "2026-02-16 04:20:45"
$ws.Range("H40").NumberFormat = "@"
$ws.Range("H40").Value = "95%"
$ws.Range("N40").Value = "2.2 °C 3:42 TU"
$ws.Range("O40").Value = "3.1 °C"
$ws.Range("E41").Value = "2026-02-16 04:20:48"
$ws.Range("H41").NumberFormat = "@"
$ws.Range("H41").Value = "51%"
$ws.Range("J41").Value = "1016.1 hPa"
$ws.Range("E42").Value = "2026-02-16 04:20:50"
$ws.Range("E43").Value = "2026-02-16 04:20:53"
$ws.Range("O43").Value = "3.5 °C"
$ws.Range("E44").Value = "2026-02-16 04:20:56"
$ws.Range("H44").NumberFormat = "@"
$ws.Range("H44").Value = "86%"
$ws.Range("O44").Value = "-0.3 °C"
$ws.Range("E45").Value = "2026-02-16 04:20:59"
$ws.Range("J45").Value = "1019.8 hPa"
$ws.Range("N45").Value = "2.2 °C 3:55 TU"
$ws.Range("O45").Value = "3.1 °C"
$ws.Range("E46").Value = "2026-02-16 04:21:01"
